# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" (strikeout) column values (column G) for rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 9
    3  = 3
    4  = 3
    5  = 1
    6  = 3
    7  = 3
    8  = 3
    9  = 5
    10 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
